# Apply the dated-worksheet update: change the date and all the
# three-digit-by-one-digit multiplication problems.

$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-04-09 Tuesday"; new = "2024-04-10 Wednesday"},
    @{old = "531×9="; new = "422×4="},
    @{old = "864×6="; new = "478×8="},
    @{old = "415×6="; new = "937×6="},
    @{old = "938×2="; new = "178×3="},
    @{old = "590×5="; new = "909×5="},
    @{old = "996×7="; new = "105×4="},
    @{old = "860×8="; new = "111×3="},
    @{old = "161×8="; new = "143×5="},
    @{old = "599×7="; new = "654×4="},
    @{old = "692×4="; new = "979×8="},
    @{old = "164×3="; new = "189×6="},
    @{old = "237×4="; new = "784×5="},
    @{old = "838×4="; new = "213×4="},
    @{old = "536×4="; new = "318×5="},
    @{old = "923×7="; new = "429×8="},
    @{old = "837×2="; new = "583×4="},
    @{old = "914×2="; new = "973×2="},
    @{old = "985×7="; new = "186×5="},
    @{old = "444×2="; new = "390×3="},
    @{old = "264×5="; new = "595×5="},
    @{old = "501×2="; new = "800×3="},
    @{old = "861×2="; new = "320×8="},
    @{old = "924×5="; new = "348×6="},
    @{old = "624×9="; new = "371×6="},
    @{old = "509×3="; new = "114×2="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
